# feat: add 2022-Q1 data
# - insert a new "2022-Q1" sheet (same layout as the other quarterly sheets)
#   right before the "总计" (totals) sheet
# - prepend a 2022-Q1 summary row to the "总计" sheet, shifting the existing
#   rows down

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet positioned right before "总计"
# ---------------------------------------------------------------------------
$wsTotalRef = $wb.Worksheets.Item(3)          # "总计" is the 3rd tab today
$ws2022 = $wb.Worksheets.Add($wsTotalRef)     # inserted immediately before it
$ws2022.Name = "2022-Q1"

# Re-fetch a fresh handle to the quarterly template sheet ("2021-Q4") and copy
# its formatting (header row + index-column styling) onto the new sheet so it
# matches the look of the other quarter tabs.
$wsQ4 = $wb.Worksheets.Item("2021-Q4")
$wsQ4.Range("A1:H13").Copy()
$ws2022.Range("A1:H13").PasteSpecial(-4122)

# Row 14 needs the same index-column styling as the other index cells
# (the template sheet only has 13 rows).
$ws2022.Range("A2").Copy()
$ws2022.Range("A14").PasteSpecial(-4122)

# Make sure the numeric-looking text columns (fund size / position / rank
# percentages) are stored as text, matching the source data.
$ws2022.Range("B2:G14").NumberFormat = "@"

# -- header row --------------------------------------------------------------
$ws2022.Cells.Item(1, 2).Value = "基金代码"
$ws2022.Cells.Item(1, 3).Value = "基金名称"
$ws2022.Cells.Item(1, 4).Value = "基金规模"
$ws2022.Cells.Item(1, 5).Value = "股票总仓位"
$ws2022.Cells.Item(1, 6).Value = "仓位占比"
$ws2022.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws2022.Cells.Item(1, 8).Value = "仓位排名"

# -- data rows ----------------------------------------------------------------
$ws2022.Cells.Item(2, 2).Value = '506005'
$ws2022.Cells.Item(2, 3).Value = '博时科创板三年定期开放混合'
$ws2022.Cells.Item(2, 4).Value = '22.84'
$ws2022.Cells.Item(2, 5).Value = '96.44'
$ws2022.Cells.Item(2, 6).Value = '3.47'
$ws2022.Cells.Item(2, 7).Value = '0.7925'
$ws2022.Cells.Item(2, 8).Value = 9
$ws2022.Cells.Item(3, 2).Value = '012650'
$ws2022.Cells.Item(3, 3).Value = '博时半导体主题混合型证券投资基金A'
$ws2022.Cells.Item(3, 4).Value = '7.18'
$ws2022.Cells.Item(3, 5).Value = '92.40'
$ws2022.Cells.Item(3, 6).Value = '6.25'
$ws2022.Cells.Item(3, 7).Value = '0.4488'
$ws2022.Cells.Item(3, 8).Value = 5
$ws2022.Cells.Item(4, 2).Value = '001048'
$ws2022.Cells.Item(4, 3).Value = '富国新兴产业股票'
$ws2022.Cells.Item(4, 4).Value = '7.99'
$ws2022.Cells.Item(4, 5).Value = '92.70'
$ws2022.Cells.Item(4, 6).Value = '4.30'
$ws2022.Cells.Item(4, 7).Value = '0.3436'
$ws2022.Cells.Item(4, 8).Value = 4
$ws2022.Cells.Item(5, 2).Value = '012651'
$ws2022.Cells.Item(5, 3).Value = '博时半导体主题混合型证券投资基金C'
$ws2022.Cells.Item(5, 4).Value = '3.14'
$ws2022.Cells.Item(5, 5).Value = '92.40'
$ws2022.Cells.Item(5, 6).Value = '6.25'
$ws2022.Cells.Item(5, 7).Value = '0.1962'
$ws2022.Cells.Item(5, 8).Value = 5
$ws2022.Cells.Item(6, 2).Value = '519606'
$ws2022.Cells.Item(6, 3).Value = '国泰金鑫股票'
$ws2022.Cells.Item(6, 4).Value = '4.25'
$ws2022.Cells.Item(6, 5).Value = '88.22'
$ws2022.Cells.Item(6, 6).Value = '4.19'
$ws2022.Cells.Item(6, 7).Value = '0.1781'
$ws2022.Cells.Item(6, 8).Value = 8
$ws2022.Cells.Item(7, 2).Value = '009057'
$ws2022.Cells.Item(7, 3).Value = '博时科技创新混合A'
$ws2022.Cells.Item(7, 4).Value = '9.59'
$ws2022.Cells.Item(7, 5).Value = '77.75'
$ws2022.Cells.Item(7, 6).Value = '1.80'
$ws2022.Cells.Item(7, 7).Value = '0.1726'
$ws2022.Cells.Item(7, 8).Value = 9
$ws2022.Cells.Item(8, 2).Value = '002181'
$ws2022.Cells.Item(8, 3).Value = '华安大安全主题灵活配置混合'
$ws2022.Cells.Item(8, 4).Value = '5.02'
$ws2022.Cells.Item(8, 5).Value = '87.28'
$ws2022.Cells.Item(8, 6).Value = '3.01'
$ws2022.Cells.Item(8, 7).Value = '0.1511'
$ws2022.Cells.Item(8, 8).Value = 10
$ws2022.Cells.Item(9, 2).Value = '009058'
$ws2022.Cells.Item(9, 3).Value = '博时科技创新混合C'
$ws2022.Cells.Item(9, 4).Value = '5.00'
$ws2022.Cells.Item(9, 5).Value = '77.75'
$ws2022.Cells.Item(9, 6).Value = '1.80'
$ws2022.Cells.Item(9, 7).Value = '0.0900'
$ws2022.Cells.Item(9, 8).Value = 9
$ws2022.Cells.Item(10, 2).Value = '012200'
$ws2022.Cells.Item(10, 3).Value = '新华鑫科技3个月滚动持有灵活配置混合型证券投资基金A'
$ws2022.Cells.Item(10, 4).Value = '2.04'
$ws2022.Cells.Item(10, 5).Value = '77.02'
$ws2022.Cells.Item(10, 6).Value = '3.42'
$ws2022.Cells.Item(10, 7).Value = '0.0698'
$ws2022.Cells.Item(10, 8).Value = 4
$ws2022.Cells.Item(11, 2).Value = '004091'
$ws2022.Cells.Item(11, 3).Value = '博时沪港深价值优选灵活配置混合A'
$ws2022.Cells.Item(11, 4).Value = '1.28'
$ws2022.Cells.Item(11, 5).Value = '58.32'
$ws2022.Cells.Item(11, 6).Value = '2.30'
$ws2022.Cells.Item(11, 7).Value = '0.0294'
$ws2022.Cells.Item(11, 8).Value = 4
$ws2022.Cells.Item(12, 2).Value = '003456'
$ws2022.Cells.Item(12, 3).Value = '信达澳银新目标灵活配置混合'
$ws2022.Cells.Item(12, 4).Value = '1.05'
$ws2022.Cells.Item(12, 5).Value = '86.04'
$ws2022.Cells.Item(12, 6).Value = '2.08'
$ws2022.Cells.Item(12, 7).Value = '0.0218'
$ws2022.Cells.Item(12, 8).Value = 4
$ws2022.Cells.Item(13, 2).Value = '012201'
$ws2022.Cells.Item(13, 3).Value = '新华鑫科技3个月滚动持有灵活配置混合型证券投资基金C'
$ws2022.Cells.Item(13, 4).Value = '0.52'
$ws2022.Cells.Item(13, 5).Value = '77.02'
$ws2022.Cells.Item(13, 6).Value = '3.42'
$ws2022.Cells.Item(13, 7).Value = '0.0178'
$ws2022.Cells.Item(13, 8).Value = 4
$ws2022.Cells.Item(14, 2).Value = '004092'
$ws2022.Cells.Item(14, 3).Value = '博时沪港深价值优选灵活配置混合C'
$ws2022.Cells.Item(14, 4).Value = '0.09'
$ws2022.Cells.Item(14, 5).Value = '58.32'
$ws2022.Cells.Item(14, 6).Value = '2.30'
$ws2022.Cells.Item(14, 7).Value = '0.0021'
$ws2022.Cells.Item(14, 8).Value = 4

# -- index column (A), numeric 0..12 ------------------------------------------
$ws2022.Cells.Item(2, 1).Value = 0
$ws2022.Cells.Item(3, 1).Value = 1
$ws2022.Cells.Item(4, 1).Value = 2
$ws2022.Cells.Item(5, 1).Value = 3
$ws2022.Cells.Item(6, 1).Value = 4
$ws2022.Cells.Item(7, 1).Value = 5
$ws2022.Cells.Item(8, 1).Value = 6
$ws2022.Cells.Item(9, 1).Value = 7
$ws2022.Cells.Item(10, 1).Value = 8
$ws2022.Cells.Item(11, 1).Value = 9
$ws2022.Cells.Item(12, 1).Value = 10
$ws2022.Cells.Item(13, 1).Value = 11
$ws2022.Cells.Item(14, 1).Value = 12

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: add a 2022-Q1 row on top, push the
#    existing 2021-Q4 / 2021-Q3 rows down by one.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q1"
$wsTotal.Cells.Item(2, 3).Value = 13
$wsTotal.Cells.Item(2, 4).Value = 2.51

$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(3, 2).Value = "2021-Q4"
$wsTotal.Cells.Item(3, 3).Value = 12
$wsTotal.Cells.Item(3, 4).Value = 2.33

$wsTotal.Cells.Item(4, 1).Value = 2
$wsTotal.Cells.Item(4, 2).Value = "2021-Q3"
$wsTotal.Cells.Item(4, 3).Value = 3
$wsTotal.Cells.Item(4, 4).Value = 0.11

# the new index cell (A4) needs the same styling as the other index cells
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A4").PasteSpecial(-4122)
